$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate row 78 (values + formatting) into a new row 79 so the new
# entry inherits the same per-column styles as the rest of the table,
# then overwrite the TestName/ID cells with the new test case.
$ws.Rows.Item(78).Copy()
$ws.Rows.Item(79).Insert()
$ws.Range("A79").Value = "Domestic_Payments_Modify_Data_[WEB]"
$ws.Range("B79").Value = "C70834"

# Re-apply the autofilter so its range grows by one row, matching the
# sheet's existing convention of leaving the very last row outside of it.
$ws.AutoFilterMode = $false
$ws.Range("A1:F77").AutoFilter()

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the
# autofilter range (Excel updates this automatically; do it explicitly
# here since the autofilter was re-created).
$fdb = $wb.Names.Item("_xlnm._FilterDatabase")
$fdb.RefersTo = '=Sheet1!$A$1:$F$78'

# Move the saved selection the same way the author's session ended up.
$ws.Range("B81").Select()
